# Update the GenBank submission TSV-style sheet:
#  - rename header columns (Sequence_ID/SRA_accession/bioproject_accession/biosample_accession
#    -> Accession/SRA/Bioproject/Biosample)
#  - replace column A "Sequence_ID" values with the new GenBank Accession numbers
#  - move the active selection to E11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Accession"
$ws.Range("B1").Value = "SRA"
$ws.Range("C1").Value = "Bioproject"
$ws.Range("D1").Value = "Biosample"

# New Accession values for column A (rows 2-30)
$accessions = @(
    "PP766449",
    "PP766450",
    "PP766451",
    "PP766452",
    "PP766453",
    "PP766454",
    "PP766455",
    "PP766456",
    "PP766457",
    "PP766458",
    "PP766459",
    "PP766461",
    "PP766460",
    "PP766462",
    "PP766463",
    "PP766464",
    "PP766465",
    "PP766466",
    "PP766467",
    "PP766468",
    "PP766469",
    "PP766470",
    "PP766471",
    "PP766472",
    "PP766473",
    "PP766474",
    "PP766475",
    "PP766476",
    "PP766477"
)

for ($i = 0; $i -lt $accessions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $accessions[$i]
}

# Update the active selection shown in the workbook
$ws.Range("E11").Select()
